$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 396
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
